$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.066.79"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.372.56"
$ws.Range("E3").Value = "  +2.25%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.12"
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.51"
$ws.Range("E6").Value = "  +1.68%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  +0.57%  "

# Row 9
$ws.Range("E9").Value = "  -2.40%  "

# Row 10
$ws.Range("E10").Value = "  +1.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("E11").Value = "  +3.05%  "

# Row 12
$ws.Range("E12").Value = "  +1.00%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.39"
$ws.Range("E13").Value = "  -1.10%  "

# Row 14
$ws.Range("E14").Value = "  +1.02%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.737.13"
$ws.Range("E15").Value = "  +1.94%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.396.23"
$ws.Range("E16").Value = "  +4.29%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  +1.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.076.56"
$ws.Range("E18").Value = "  +0.88%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.00"
$ws.Range("E19").Value = "  -0.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.26"
$ws.Range("E20").Value = "  +1.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  +0.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.98"
$ws.Range("E22").Value = "  +0.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.50"
$ws.Range("E23").Value = "  +0.22%  "

# Row 24
$ws.Range("E24").Value = "  -1.14%  "

# Row 25
$ws.Range("E25").Value = "  +0.96%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.44"
$ws.Range("E27").Value = "  +0.38%  "

# Row 28
$ws.Range("E28").Value = "  +0.76%  "

# Row 29
$ws.Range("E29").Value = "  +2.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.34"
$ws.Range("E30").Value = "  +3.95%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.06%  "

# Row 32
$ws.Range("E32").Value = "  +1.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.61"
$ws.Range("E33").Value = "  +2.92%  "

# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  +9.40%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0730"
$ws.Range("E35").Value = "  -3.24%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  +1.95%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "127.07"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38
$ws.Range("E38").Value = "  +5.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.32"
$ws.Range("E39").Value = "  -0.61%  "

# Row 40
$ws.Range("E40").Value = "  -2.25%  "

# Row 41
$ws.Range("E41").Value = "  -0.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.82"
$ws.Range("E42").Value = "  -5.18%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.929.06"
$ws.Range("E43").Value = "  +0.39%  "

# Row 44
$ws.Range("E44").Value = "  -0.42%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("E45").Value = "  +2.64%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  +1.15%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.21"
$ws.Range("E47").Value = "  -8.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.597.54"
$ws.Range("E48").Value = "  +1.77%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.50"
$ws.Range("E49").Value = "  +3.46%  "

# Row 50
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("E51").Value = "  +1.26%  "
